# Reassign duty-activity labels across the four timetable sheets.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: CSE-IV (B)
# ---------------------------------------------------------------
$wsB = $wb.Worksheets.Item("CSE-IV (B)")

$wsB.Range("H10").Value = "PROJECT WORK"

$wsB.Range("C11:D11").UnMerge()
$wsB.Range("B11").Value = "LIB"
$wsB.Range("C11").Value = "OB"
$wsB.Range("D11").Value = "HCI"
$wsB.Range("H11").Value = "PROJECT WORK"

$wsB.Range("H10:J10").Merge()
$wsB.Range("D11:E11").Merge()
$wsB.Range("H11:J11").Merge()

# ---------------------------------------------------------------
# Sheet: CSE-IV (C)
# ---------------------------------------------------------------
$wsC = $wb.Worksheets.Item("CSE-IV (C)")

$wsC.Range("B10").Value = "LIB"
$wsC.Range("I10").Value = "SEMINAR"
$wsC.Range("G11").Value = "OB"
$wsC.Range("H11").Value = "HCI"

$wsC.Range("I10:J10").Merge()

# ---------------------------------------------------------------
# Sheet: CSE-IV (D)
# ---------------------------------------------------------------
$wsD = $wb.Worksheets.Item("CSE-IV (D)")

$wsD.Range("B10").Value = "LIB"
$wsD.Range("H10").Value = "PROJECT WORK"
$wsD.Range("B11").Value = "LIB"
$wsD.Range("G11").Value = "SEMINAR"

$wsD.Range("H10:J10").Merge()

# ---------------------------------------------------------------
# Sheet: IT-IV
# ---------------------------------------------------------------
$wsIT = $wb.Worksheets.Item("IT-IV")

$wsIT.Range("H10:I10").UnMerge()
$wsIT.Range("B10").Value = "PROJECT WORK"
$wsIT.Range("G10").Value = "HCI"
$wsIT.Range("H10").Value = ""
$wsIT.Range("I10").Value = "OB"
$wsIT.Range("J10").Value = "SPORTS"
$wsIT.Range("B11").Value = "SEMINAR"

$wsIT.Range("B10:D10").Merge()
$wsIT.Range("G10:H10").Merge()
$wsIT.Range("B11:C11").Merge()

# Keep IT-IV as the active/selected sheet & cell, matching the source file.
$wsIT.Activate()
$wsIT.Range("I11").Select()
